$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New lookup column (C) next to the existing gear-type list, plus a
# --- small "OLD GEAR" cross-reference table lower on the sheet (gear-name
# --- clean-up / CDOM data processing pass).

$ws.Range("C1").Value = "replace existin gear in database"

$ws.Range("C3").Value = "Water bottle 10l"
$ws.Range("C4").Value = "bucket; sieve/Jar; Surface"
$ws.Range("C6").Value = "20 μm net; 20μm net; Phytoplankton net"
$ws.Range("C9").Value = "Ice core"
$ws.Range("C13").Value = "MIK"
$ws.Range("C14").Value = "Multinet; MPS"
$ws.Range("C16").Value = "Niskin"
$ws.Range("C18").Value = "CTD"
$ws.Range("C22").Value = "suction pump; suction pump & net"
$ws.Range("C23").Value = "hand net; Square net"
$ws.Range("C26").Value = "WP2 150"
$ws.Range("C27").Value = "WP2 180; WP2 200"
$ws.Range("C28").Value = "WP2 500"
$ws.Range("C29").Value = "WP2 63; WP2 63μm; WP2 60um"
$ws.Range("C30").Value = "WP3 1000"

# New gear-type row
$ws.Range("A31").Value = "Diver"
$ws.Range("C31").Value = "Divers; Divers; Snorkling"

# "OLD GEAR" mini table
$ws.Range("C33").Value = "OLD GEAR"

$ws.Range("C34").Value = "plankton net"
$ws.Range("D34").Value = "check wich sample type, could be Handnet 20 µm"

$ws.Range("C35").Value = "Mega zooplankton net 1.55"
$ws.Range("D35").Value = "same as `"Macro zooplankton net`"; keep one of the names"

$ws.Range("C36").Value = "juday net"
$ws.Range("D36").Value = "delete if no sample is attached"

$ws.Range("C37").Value = "Grab"
$ws.Range("D37").Value = "delete if no sample is attached"

$ws.Range("C38").Value = "Bongo"
$ws.Range("D38").Value = "delete if no sample is attached"

$ws.Range("C39").Value = "RP sledge"
$ws.Range("D39").Value = "delete if no sample is attached"

# --- Formatting: reuse existing cell styles instead of inventing new ones.

# Bold header style (same as B1) -> C1 and C33
$boldHeaderCells = @("C1", "C33")
foreach ($addr in $boldHeaderCells) {
    $ws.Range("B1").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

# Plain "Calibri 10, theme text" style (same as A7) -> every other new C/D cell
$plainCells = @("C3","C4","C6","C9","C13","C14","C16","C18","C22","C23",
                "C26","C27","C28","C29","C30","C31")
foreach ($addr in $plainCells) {
    $ws.Range("A7").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}
$ws.Range("A7").Copy()
$ws.Range("C34:D39").PasteSpecial(-4122)

# A31 ("Diver") should look like the rest of column A (same as A30)
$ws.Range("A30").Copy()
$ws.Range("A31").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Rich-text run inside D34: "...Handnet 20 " stays default, "µm" gets an
# --- explicit (non-themed) Calibri 10 run, matching the author's manual
# --- partial-selection font tweak.
$chars = $ws.Range("D34").Characters(45, 2)
$chars.Font.Size = 10
$chars.Font.Name = "Calibri"

# --- Column widths for the two new columns (characters).
$ws.Columns.Item(3).ColumnWidth = 32.71
$ws.Columns.Item(4).ColumnWidth = 46.85

# --- Selection left where the author's cursor ended up.
$ws.Range("A9").Select() | Out-Null
